# Rebuild the "Ventas" sheet: replace the old id_venta / id_producto /
# cantidad / nombreCliente / fecha (A1:F2) layout with a new
# id_producto / cantidad / nombreCliente table (A1:C10, 9 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture the existing header format (style index 1: bold, bordered,
# centered/top-aligned) from B1 *before* clearing, so we can stamp it back
# onto the rebuilt header without minting a brand-new style entry.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> stash on Z1

# Wipe only the old used range (A1:F2) - leave the Z1 stash alone.
$ws.Range("A1:F2").Clear()

# --- header row -------------------------------------------------------
$ws.Range("A1").Value = "id_producto"
$ws.Range("B1").Value = "cantidad"
$ws.Range("C1").Value = "nombreCliente"

$ws.Range("Z1").Copy() | Out-Null
$ws.Range("A1:C1").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Clear()

# --- data rows ----------------------------------------------------------
$data = @(
    @("673d509143a0b377f4f6de02", 5,  "Juan Pérez"),
    @("673d509143a0b377f4f6de03", 3,  "María López"),
    @("673d509143a0b377f4f6de04", 12, "Carlos García"),
    @("673d509143a0b377f4f6de05", 2,  "Ana Fernández"),
    @("673d509143a0b377f4f6de06", 6,  "Luis Martínez"),
    @("673d509143a0b377f4f6de07", 10, "Sofía Gómez"),
    @("673d509143a0b377f4f6de08", 8,  "Martín Torres"),
    @("673d509143a0b377f4f6de09", 4,  "Clara Suárez"),
    @("673d509143a0b377f4f6de10", 15, "Diego Vargas")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# --- body style: vertical-center + wrap text, built once on A2, then
# fanned out to the rest of the body via copy/paste-special so only a
# single new cellXf is minted.
$ws.Range("A2").VerticalAlignment = -4108   # xlCenter
$ws.Range("A2").WrapText = $true

$ws.Range("A2").Copy() | Out-Null
$ws.Range("A2:C10").PasteSpecial(-4122) | Out-Null

# --- column widths (bestFit in the target file; input chosen so the
# engine's 1/6-character quantization lands on the closest achievable
# width to the target's 25.42578125 / 8.5703125 / 15.85546875) ---------
$ws.Columns.Item(1).ColumnWidth = 24.666666666666668
$ws.Columns.Item(2).ColumnWidth = 7.666666666666667
$ws.Columns.Item(3).ColumnWidth = 15

# --- selection matches the target sheetView -----------------------------
$ws.Range("E10").Select()

$excel.CutCopyMode = $false
